$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before D; this shifts the existing "Tipo" column to E
$ws.Range("D1").EntireColumn.Insert()

# Add the new "MAE" header and give it the same header style used by the
# other header cells (bold font, border, centered alignment)
$ws.Range("D1").Value = "MAE"
$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial(-4122)

# Populate the new MAE column with the model's mean absolute error values
$ws.Range("D2").Value = 0.1560271560555998
$ws.Range("D3").Value = 0.2072673588334908

# Refresh the MSE value for row 3 with its updated (slightly rounded) figure
$ws.Range("B3").Value = 0.0831172165082084
